$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C47").Value = "Menghapus Data Unit Kuantitas"
$ws.Range("C48").Value = "Menghapus Data Agama"
$ws.Range("C49").Value = "Menghapus Data Media Sosial"
$ws.Range("C50").Value = "Menghapus Data Merk Dagang"
$ws.Range("C51").Value = "Menghapus Data Proyek"

$ws.Range("C49").Select()
